# Edit workbook: populate missing metadata columns on the 保險 (Insurance)
# and 事業投資 (Business investment) sheets, and rename the "otherbonds"
# category label to "antique" (used by the 具有相當價值之財產 sheet).
#
# Commit message: "#5: insurance, claim, debt, investment done"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the shared "otherbonds" category string to "antique".
#    This string is used by the "具有相當價值之財產" sheet (sheet order 5)
#    and simply gets corrected/renamed in place.
# ---------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("otherbonds", "antique")
}

# ---------------------------------------------------------------------
# 2) 保險 (Insurance) sheet - 6th worksheet
#    Fix the header row (it had stray data values instead of column
#    names) and append the standard metadata columns E:K to every row.
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Header row
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

$insLastRow = 8
$wsIns.Range("G2:G8").NumberFormat = "@"
for ($r = 2; $r -le $insLastRow; $r++) {
    $idx = $wsIns.Cells.Item($r, 1).Value2
    $wsIns.Cells.Item($r, 5).Value = "insurance"
    $wsIns.Cells.Item($r, 6).Value = "normal"
    $wsIns.Cells.Item($r, 7).Value = "2012-04-19"
    $wsIns.Cells.Item($r, 8).Value = "張慶忠"
    $wsIns.Cells.Item($r, 9).Value = 1347
    $wsIns.Cells.Item($r, 10).Value = "tmp93201"
    $wsIns.Cells.Item($r, 11).Value = $idx
}

# ---------------------------------------------------------------------
# 3) 事業投資 (Business investment) sheet - 7th worksheet
#    Fix the header row and append the standard metadata columns H:N.
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("事業投資")

# Header row
$wsInv.Range("B1").Value = "owner"
$wsInv.Range("C1").Value = "company"
$wsInv.Range("D1").Value = "address"
$wsInv.Range("E1").Value = "total"
$wsInv.Range("F1").Value = "register_date"
$wsInv.Range("G1").Value = "register_reason"
$wsInv.Range("H1").Value = "property_category"
$wsInv.Range("I1").Value = "category"
$wsInv.Range("J1").Value = "date"
$wsInv.Range("K1").Value = "legislator_name"
$wsInv.Range("L1").Value = "legislator_id"
$wsInv.Range("M1").Value = "source_file"
$wsInv.Range("N1").Value = "index"

$invLastRow = 5
$wsInv.Range("J2:J5").NumberFormat = "@"
for ($r = 2; $r -le $invLastRow; $r++) {
    $idx = $wsInv.Cells.Item($r, 1).Value2
    $wsInv.Cells.Item($r, 8).Value = "investment"
    $wsInv.Cells.Item($r, 9).Value = "normal"
    $wsInv.Cells.Item($r, 10).Value = "2012-04-19"
    $wsInv.Cells.Item($r, 11).Value = "張慶忠"
    $wsInv.Cells.Item($r, 12).Value = 1347
    $wsInv.Cells.Item($r, 13).Value = "tmp93201"
    $wsInv.Cells.Item($r, 14).Value = $idx
}
